$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 becomes what used to be row 7's data (Resnet w/Quads)
$ws.Range("A6").Value = "Resnet w/Quads"
$ws.Range("B6").Value = 80
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = 34
$ws.Range("E6").Value = 39
$ws.Range("F6").Value = 12
$ws.Range("G6").Value = 46
$ws.Range("H6").Formula = "=AVERAGE(C6:G6)"

# Row 7 is new data with no label (old A7/B7 content is gone)
$ws.Range("A7").Clear()
$ws.Range("B7").Clear()
$ws.Range("C7").Value = 154
$ws.Range("D7").Value = 185
$ws.Range("E7").Value = 144
$ws.Range("F7").Value = 200
$ws.Range("G7").Value = 169
$ws.Range("H7").Formula = "=AVERAGE(C7:G7)"

# Row 11 gets what used to be row 6's data (Deeper w/Quads PG300)
$ws.Range("A11").Value = "Deeper w/Quads PG300"
$ws.Range("B11").Value = 78
$ws.Range("B11").HorizontalAlignment = -4131
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 52
$ws.Range("E11").Value = 33
$ws.Range("F11").Value = 25
$ws.Range("G11").Value = 0
$ws.Range("H11").Formula = "=AVERAGE(C11:G11)"

# Update selection
$ws.Range("H6:H7").Select()

# Page setup orientation change (1 = xlPortrait)
$ws.PageSetup.Orientation = 1
